# "new version of diploma" - populate the pressure/voltage (A/B) columns for
# the second data table (rows 67-83 on Лист1) and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Values for column A (давление, дел) and column B, rows 67-82
$dataA = @{
    67 = 3300
    68 = 2600
    69 = 2600
    70 = 1600
    71 = 1200
    72 = 800
    73 = 400
    74 = 8.954
    76 = 4000
    77 = 3600
    78 = 2300
    79 = 1800
    80 = 1400
    81 = 1000
    82 = 800
    83 = 8.317
}

$dataB = @{
    67 = 50
    68 = 250
    69 = 450
    70 = 650
    71 = 850
    72 = 1050
    73 = 1250
    74 = 0.0023
    76 = 50
    77 = 650
    78 = 1250
    79 = 1850
    80 = 2450
    81 = 3050
    82 = 3650
    83 = 0.0004
}

foreach ($r in 67..83) {
    if ($dataA.ContainsKey($r)) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.Value = $dataA[$r]
        $cell.HorizontalAlignment = -4108
    }
    if ($dataB.ContainsKey($r)) {
        $cell = $ws.Cells.Item($r, 2)
        $cell.Value = $dataB[$r]
        $cell.HorizontalAlignment = -4108
    }
}

# Move the active selection to C83, matching the post-edit sheet view.
$ws.Range("C83").Select()
